# Adds four new application-register rows (49-52) describing new
# self-pickup ("самовывоз") shipments, including the organization,
# transshipment point and purchaser details, plus the free-text request.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 49
$ws.Range("B49").Value = "самовывоз"
$ws.Range("C49").Value = "28.04.2024"
$ws.Range("D49").Value = "Цем I 42,5Н/ Кричев "
$ws.Range("F49").Value = 70
$ws.Range("G49").Value = "т"
$ws.Range("H49").Value = "н223вм977 Музафяров, н095вм977 Трушев, Т132ак977 Бабий, Т132Ак977 Савинец, Н223ВМ977 Маркин, Н095ВМ977 Шаркевич"
$ws.Range("K49").Value = "`"ЗАО ПК ТЕРМОБЕТОН  "
$ws.Range("R49").Value = " Заявка на самовывоз  1. Дата отгрузки:  28.04.2024  Марка Цем I 42,5Н/ Кричев  3. Количество машин/тонн: 70 тонн  4. Перевалка Сзтк  5.Покупатель груза: `"ЗАО ПК ТЕРМОБЕТОН   6. Продажа от ООО «СЗТК»  н223вм977 Музафяров Руслан н095вм977 Трушев Давид Т132ак977 Бабий Антон Т132Ак977 Савинец Юрий Н223ВМ977 Маркин Александр Н095ВМ977 Шаркевич Алексей "

# Row 50
$ws.Range("A50").Value = "Юрий"
$ws.Range("B50").Value = "самовывоз"
$ws.Range("C50").Value = "26.04.2024"
$ws.Range("D50").Value = ": ЦЕМ I 42.5Н БЦЗ Костюковичи "
$ws.Range("F50").Value = 100
$ws.Range("G50").Value = "т"
$ws.Range("K50").Value = "ООО `"Форма ЖБИ`""
$ws.Range("R50").Value = " Юра Менеджер: Заявка/ самовывоз 1. Дата отгрузки: 26.04.24 2. Марка цемента: ЦЕМ I 42.5Н БЦЗ Костюковичи  3. Кол-во машин/ тонн: 100т 4. Продажа от клиента: ОО0 Спарта  5. Завод отгрузки: СЗТК 6. Покупатель груза:ООО `"Форма ЖБИ`" 7. Грузополучатель (при оформлении ттн): ООО `"Форма ЖбИ`" 8. М192МС90"

# Row 51
$ws.Range("A51").Value = "Юрий"
$ws.Range("B51").Value = "самовывоз"
$ws.Range("C51").Value = "26.04.2024"
$ws.Range("D51").Value = ": ЦЕМ I 42.5Н БЦЗ Костюковичи "
$ws.Range("F51").Value = 100
$ws.Range("G51").Value = "т"
$ws.Range("I51").Value = "ОО0 Спарта "
$ws.Range("K51").Value = "ООО `"Форма ЖБИ`""
$ws.Range("R51").Value = " Юра Менеджер: Заявка/ самовывоз 1. Дата отгрузки: 26.04.24 2. Марка цемента: ЦЕМ I 42.5Н БЦЗ Костюковичи  3. Кол-во машин/ тонн: 100т 4. Продажа от клиента: ОО0 Спарта  5. Завод отгрузки: СЗТК 6. Покупатель груза:ООО `"Форма ЖБИ`" 7. Грузополучатель (при оформлении ттн): ООО `"Форма ЖбИ`" 8. М192МС90"

# Row 52
$ws.Range("A52").Value = "Юрий"
$ws.Range("B52").Value = "самовывоз"
$ws.Range("C52").Value = "26.04.2024"
$ws.Range("D52").Value = ": ЦЕМ I 42.5Н БЦЗ Костюковичи "
$ws.Range("F52").Value = 100
$ws.Range("G52").Value = "т"
$ws.Range("I52").Value = "ОО0 Спарта "
$ws.Range("J52").Value = "СЗТК"
$ws.Range("K52").Value = "ООО `"Форма ЖБИ`""
$ws.Range("R52").Value = " Юра Менеджер: Заявка/ самовывоз 1. Дата отгрузки: 26.04.24 2. Марка цемента: ЦЕМ I 42.5Н БЦЗ Костюковичи  3. Кол-во машин/ тонн: 100т 4. Продажа от клиента: ОО0 Спарта  5. Завод отгрузки: СЗТК 6. Покупатель груза:ООО `"Форма ЖБИ`" 7. Грузополучатель (при оформлении ттн): ООО `"Форма ЖбИ`" 8. М192МС90"
